# Scheduled runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a handful of leve
# rows on the ALC, CRP, CUL and LTW sheets. All cells here are plain
# cached numbers (no formulas), so this just overwrites each changed
# cell with its new value; a few cells that become genuinely empty are
# cleared instead of zero-filled.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1433.8975
$ws.Range("I98").Value = 512.29034
$ws.Range("J98").Value = 5005.125
$ws.Range("K98").Value = 512.29034
$ws.Range("L98").Value = 5005.125
$ws.Range("M98").Value = 985.70966
$ws.Range("N98").Value = -8001.125
$ws.Range("H107").Value = 312.42105
$ws.Range("I107").Value = 259
$ws.Range("J107").Value = 512.75
$ws.Range("K107").Value = 259
$ws.Range("L107").Value = 512.75
$ws.Range("M107").Value = 1661
$ws.Range("N107").Value = -4352.75
$ws.Range("H122").Value = 1433.8975
$ws.Range("I122").Value = 512.29034
$ws.Range("J122").Value = 5005.125
$ws.Range("K122").Value = 1536.87102
$ws.Range("L122").Value = 15015.375
$ws.Range("M122").Value = 913.12898
$ws.Range("N122").Value = -19915.375
$ws.Range("H137").Value = 5655.615
$ws.Range("I137").Value = 6862.4
$ws.Range("J137").Value = 1633
$ws.Range("K137").Value = 20587.2
$ws.Range("L137").Value = 4899
$ws.Range("M137").Value = -18037.2
$ws.Range("N137").Value = -9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43798.805
$ws.Range("I31").Value = 51311.92
$ws.Range("K31").Value = 51311.92
$ws.Range("M31").Value = -51016.92
$ws.Range("H34").Value = 43798.805
$ws.Range("I34").Value = 51311.92
$ws.Range("K34").Value = 51311.92
$ws.Range("M34").Value = -51109.92

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1869.925
$ws.Range("I5").Value = 672
$ws.Range("K5").Value = 2016
$ws.Range("M5").Value = -1904
$ws.Range("H62").Value = 3140
$ws.Range("I62").Value = 1300
$ws.Range("J62").Value = 4980
$ws.Range("K62").Value = 3900
$ws.Range("L62").Value = 14940
$ws.Range("M62").Value = -3214
$ws.Range("N62").Value = -16312
$ws.Range("H63").Value = 4155
$ws.Range("I63").Value = 3732.5
$ws.Range("K63").Value = 11197.5
$ws.Range("M63").Value = -10448.5
$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 12000
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = -11730
$ws.Range("N64").Value = -12540
$ws.Range("H65").Value = 3140
$ws.Range("I65").Value = 1300
$ws.Range("J65").Value = 4980
$ws.Range("K65").Value = 11700
$ws.Range("L65").Value = 44820
$ws.Range("M65").Value = -8268
$ws.Range("N65").Value = -51684
$ws.Range("H66").Value = 4155
$ws.Range("I66").Value = 3732.5
$ws.Range("K66").Value = 33592.5
$ws.Range("M66").Value = -29848.5
$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 12000
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = -11064
$ws.Range("N67").Value = -13872
$ws.Range("H70").Value = 2202
$ws.Range("I70").Value = 876.3
$ws.Range("J70").Value = 4411.5
$ws.Range("K70").Value = 2628.9
$ws.Range("L70").Value = 13234.5
$ws.Range("M70").Value = -2313.9
$ws.Range("N70").Value = -13864.5
$ws.Range("H73").Value = 2202
$ws.Range("I73").Value = 876.3
$ws.Range("J73").Value = 4411.5
$ws.Range("K73").Value = 2628.9
$ws.Range("L73").Value = 13234.5
$ws.Range("M73").Value = -1536.9
$ws.Range("N73").Value = -15418.5
$ws.Range("H80").Value = 4984
$ws.Range("I80").Value = 4980
$ws.Range("J80").Value = 4985
$ws.Range("K80").Value = 14940
$ws.Range("L80").Value = 14955
$ws.Range("M80").Value = -14004
$ws.Range("N80").Value = -16827
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H82").Value = 18000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 18000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 54000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -54812
$ws.Range("H83").Value = 4984
$ws.Range("I83").Value = 4980
$ws.Range("J83").Value = 4985
$ws.Range("K83").Value = 44820
$ws.Range("L83").Value = 44865
$ws.Range("M83").Value = -40140
$ws.Range("N83").Value = -54225
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H85").Value = 18000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 18000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 54000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -56808
$ws.Range("H95").Value = 5000
$ws.Range("I95").Value = 5000
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 15000
$ws.Range("L95").Value = 0
$ws.Range("M95").Value = -12941
$ws.Range("N95").ClearContents()
$ws.Range("H96").Value = 4520
$ws.Range("I96").Value = 3200
$ws.Range("J96").Value = 4850
$ws.Range("K96").Value = 9600
$ws.Range("L96").Value = 14550
$ws.Range("M96").Value = -7541
$ws.Range("N96").Value = -18668
$ws.Range("H97").Value = 751.5454999999999
$ws.Range("I97").Value = 430.42856
$ws.Range("J97").Value = 1313.5
$ws.Range("K97").Value = 1291.28568
$ws.Range("L97").Value = 3940.5
$ws.Range("M97").Value = -795.28568
$ws.Range("N97").Value = -4932.5
$ws.Range("H132").Value = 1538.3846
$ws.Range("I132").Value = 949.8333
$ws.Range("J132").Value = 2042.8572
$ws.Range("K132").Value = 8548.4997
$ws.Range("L132").Value = 18385.7148
$ws.Range("M132").Value = -6018.4997
$ws.Range("N132").Value = -23445.7148
$ws.Range("H135").Value = 1869.925
$ws.Range("I135").Value = 672
$ws.Range("K135").Value = 6048
$ws.Range("M135").Value = -3513

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8606
$ws.Range("I132").Value = 10250
$ws.Range("J132").Value = 5318
$ws.Range("K132").Value = 30750
$ws.Range("L132").Value = 15954
$ws.Range("M132").Value = -28220
$ws.Range("N132").Value = -21014
